# Insert a new slide "Используемые модули" (Used modules) right after the
# title slide (new slide #2); every slide that followed shifts down by one.
# This mirrors the commit's sldIdLst change: a brand new slide id (267) is
# inserted between id 256 (slide 1) and id 257 (old slide 2).

$p = $ppt.ActivePresentation

# ppLayoutText = 2 -> "Title and Content" layout (same "Заголовок и объект"
# layout, slideLayout2.xml, used by every other content slide in the deck).
$s = $p.Slides.Add(2, 2)

# The new slide only ends up with a body placeholder (idx=1) and a free
# floating title textbox in the final deck -- drop the auto-created Title
# placeholder shape entirely.
$s.Shapes.Item(1).Delete()

# Reuse the body/content placeholder that Slides.Add already created
# (keeps it a real placeholder, p:ph idx="1", matching the target XML).
$body = $s.Shapes.Item(1)
$body.Name = "Объект 2"

# Position/size in points (COM reports Left/Top/Width/Height in points,
# 1 pt = 12700 EMU) so the saved EMU box matches 685800,2096906,8156359,3433883.
$body.Left = 54.0
$body.Top = 165.11071015141852
$body.Width = 642.2330017159918
$body.Height = 270.3844881889764

$bodyText = $body.TextFrame.TextRange
$bodyText.Text = "selenium`rfuzzywuzzy`rsympy`rdiscord`rPyttsx3`ryoutube_dl"

# Uniform run formatting for every line: 32pt Times New Roman (latin + cs).
for ($i = 1; $i -le 6; $i++) {
    $run = $bodyText.Paragraphs($i, 1)
    $run.Font.Size = 32
    $run.Font.Name = "Times New Roman"
    $run.Font.NameComplexScript = "Times New Roman"
}
$bodyText.LanguageID = "en-US"

# Consume shape id 4 the same way the source deck does (a shape was
# created and removed before the final title textbox landed on id 5,
# named "TextBox 4") so the id/name numbering lines up with the target.
$placeholderShape = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$placeholderShape.Delete()

# Free-floating title textbox ("Используемые модули"), 60pt Times New Roman.
$title = $s.Shapes.AddTextbox(1, 54.0, 69.91260147519493, 609.3786614173229, 79.97346456692914)
$titleText = $title.TextFrame.TextRange
$titleText.Text = "Используемые модули"
$titleText.Font.Size = 60
$titleText.Font.Name = "Times New Roman"
$titleText.Font.NameComplexScript = "Times New Roman"
$titleText.LanguageID = "ru-RU"
